$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 323
$ws.Range("F3").Value = 528
$ws.Range("F4").Value = 590
$ws.Range("F5").Value = 9257
$ws.Range("F6").Value = 26
$ws.Range("F10").Value = 309
$ws.Range("F14").Value = 18
$ws.Range("F17").Value = 2069
$ws.Range("F18").Value = 828
$ws.Range("F19").Value = 785
$ws.Range("F21").Value = 45
$ws.Range("F22").Value = 411
$ws.Range("F24").Value = 98
$ws.Range("F25").Value = 665
$ws.Range("F26").Value = 20
$ws.Range("F27").Value = 1557
$ws.Range("F28").Value = 39
$ws.Range("F30").Value = 20
$ws.Range("F33").Value = 1444
$ws.Range("F34").Value = 7
$ws.Range("F35").Value = 504
$ws.Range("F37").Value = 543
$ws.Range("F38").Value = 398
$ws.Range("F39").Value = 2196
$ws.Range("F41").Value = 77
$ws.Range("F42").Value = 157
$ws.Range("F43").Value = 577
$ws.Range("F45").Value = 165
$ws.Range("F46").Value = 887
$ws.Range("F47").Value = 682
$ws.Range("F49").Value = 316
$ws.Range("F50").Value = 283

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 43
$ws.Range("F15").Value = 72
$ws.Range("F18").Value = 138
$ws.Range("F23").Value = 110
$ws.Range("F24").Value = 70

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2908
$ws.Range("F6").Value = 268

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 323
$ws.Range("F4").Value = 528
$ws.Range("F5").Value = 43
$ws.Range("F8").Value = 590
$ws.Range("F9").Value = 9257
$ws.Range("F10").Value = 26
$ws.Range("F11").Value = 11981
$ws.Range("F15").Value = 18
$ws.Range("F17").Value = 45
$ws.Range("F18").Value = 411
$ws.Range("F21").Value = 665
$ws.Range("F22").Value = 20
$ws.Range("F23").Value = 268
$ws.Range("F24").Value = 1557
$ws.Range("F25").Value = 39
$ws.Range("F31").Value = 72
$ws.Range("F32").Value = 1444
$ws.Range("F34").Value = 7
$ws.Range("F35").Value = 504
$ws.Range("F36").Value = 543
$ws.Range("F37").Value = 398
$ws.Range("F39").Value = 2196
$ws.Range("F40").Value = 77
$ws.Range("F41").Value = 157
$ws.Range("F42").Value = 577
$ws.Range("F44").Value = 165
$ws.Range("F45").Value = 887
$ws.Range("F49").Value = 682
